$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper pattern used throughout: inserting new text at a point whose
# neighbouring run(s) share identical run formatting normally gets folded
# straight back into that neighbour (the engine keeps adjacent,
# identically-formatted runs coalesced). To force a genuine run break -
# exactly like Word leaves behind after a real edit - we briefly flip a
# formatting flag (Bold) on the freshly inserted text and then flip it
# back off; that round-trip is enough to keep it as its own <w:r> even
# though the final rendered formatting matches its neighbours again.
# -----------------------------------------------------------------------

# ---------------------------------------------------------------------
# Edit 1: "...split out users based on sex." -> "...based on gender."
# and move the "_GoBack" bookmark so it sits right after "gender".
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("split out users based on sex.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "edit1 found=$found1"
if ($found1) {
    $end1 = $r1.End
    $sexRange = $d.Range($end1 - 4, $end1 - 1)
    Write-Output "sexRange text=[$($sexRange.Text)]"
    $sexRange.Text = "gender"
    $sexRange.Font.Bold = $true
    $sexRange.Font.Bold = $false

    # relocate the "_GoBack" bookmark to sit right after "gender"
    $oldBm1 = $d.Bookmarks.Item("_GoBack")
    $oldBm1.Delete()
    $bmPos1 = $d.Range($sexRange.End, $sexRange.End)
    $d.Bookmarks.Add("_GoBack", $bmPos1)
}

# ---------------------------------------------------------------------
# Edit 2: merge "orrelation Between Personality Facets" + ":" runs
# (re-running Find/Replace with identical text forces the engine to
# normalize/merge the two adjacent, identically formatted runs)
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("orrelation Between Personality Facets:", $true, $false, $false, $false, $false, $true, 1, $false, "orrelation Between Personality Facets:", 2)
Write-Output "edit2 found=$found2"

# ---------------------------------------------------------------------
# Edit 3: "...implemented OOP to create a system..." ->
#         "...implemented OOP approach to create a system..."
# ---------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("implemented OOP to create", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "edit3 found=$found3"
if ($found3) {
    $start3 = $r3.Start
    $insertAt3 = $start3 + 15
    $ip3 = $d.Range($insertAt3, $insertAt3)
    $ip3.InsertAfter(" approach")
    $ip3.Font.Bold = $true
    $ip3.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Edit 4: remove the old bookmarkStart/End (_GoBack) near
# "After identify|ing" -- already relocated in Edit 1 above, so this spot
# no longer carries the bookmark once Edit 1 has run.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Edit 5: "...implement the traversal into an online..." ->
#         "...implement the tree traversal into an online..."
# ---------------------------------------------------------------------
$r5 = $d.Content
$found5 = $r5.Find.Execute("implement the traversal into an online", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "edit5 found=$found5"
if ($found5) {
    $start5 = $r5.Start
    $insertAt5 = $start5 + 14
    $ip5 = $d.Range($insertAt5, $insertAt5)
    $ip5.InsertAfter("tree ")
    $ip5.Font.Bold = $true
    $ip5.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Edit 6: remove the empty paragraph (ind left=360, empty) that follows
# the "Once finished ... portfolio. " paragraph
# ---------------------------------------------------------------------
$targetIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "my online digital portfolio") {
        $targetIdx = $i
    }
    $i = $i + 1
}
Write-Output "edit6 targetIdx=$targetIdx"
if ($targetIdx -ge 0) {
    $j = 0
    foreach ($p in $d.Paragraphs) {
        if ($j -eq ($targetIdx + 1)) {
            Write-Output "nextPara text=[$($p.Range.Text)]"
            if ($p.Range.Text.Trim().Length -eq 0) {
                $p.Range.Delete()
            }
        }
        $j = $j + 1
    }
}

Write-Output "DONE"
